$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Date column (B) to stay plain text rather than being
# auto-converted by Excel into a date serial number, then reset the
# style back to Normal so no stray per-cell formatting is introduced.
$ws.Range("B2:B12").NumberFormat = "@"

$data = @(
  ,@("Japanese J League 3", "2025-11-10", "03:00:00", "FC Osaka", "Nagano Parceiro", 1.02, 1.68, 2.46, 1000, 2.5, 1000, 1.01, 1.01, 1.3, 1.01, 1.3, 1.32, 1.18, 1.32, 1.01, 1.01, 1.01, 2.46, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
  ,@("Romanian Liga II", "2025-11-10", "14:30:00", "ASA Targu Mures", "Poli Iasi", 1.04, 1000, 1.04, 1000, 1.02, 1000, 1.01, 1.01, 1.01, 1.02, 1.24, 1.39, 1.15, 1.39, 1.01, 1.01, 1.01, 1.01, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
  ,@("Danish 1st Division", "2025-11-10", "15:00:00", "Esbjerg", "Middelfart", 1.72, 1.88, 4.6, 5.5, 3.9, 5, 1.24, 1.04, 2.12, 1.25, 2.12, 1.73, 1.36, 2.56, 1.51, 1.78, 1.22, 2.12, 25, 29, 55, 1000, 14.5, 13.5, 29, 90, 17, 15, 28, 90, 28, 26, 48, 1000, 1000, 1000)
  ,@("English Premier League 2 - Div 1", "2025-11-10", "16:00:00", "Blackburn U21", "Derby U21", 2.12, 2.74, 2.34, 3.45, 3.8, 10, 1.01, 1.02, 1.01, 1.01, 1.24, 1.01, 1.83, 1.9, 1.01, 1.01, 1.01, 1.01, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
  ,@("Italian Serie C", "2025-11-10", "16:30:00", "Salernitana", "Crotone", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1.24, 1.01, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  ,@("Spanish Segunda Division", "2025-11-10", "16:30:00", "Burgos", "CD Castellon", 2.6, 2.86, 2.88, 3.2, 3.15, 3.4, 1.01, 1.01, 1.73, 1.01, 1.73, 1.97, 1.03, 1.97, 1.01, 1.01, 1.45, 1.54, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
  ,@("Argentinian Primera Division", "2025-11-10", "17:00:00", "Gimnasia La Plata", "Velez Sarsfield", 3.55, 3.95, 2.44, 2.66, 2.8, 3.1, 0, 0, 0, 0, 1.42, 3, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  ,@("Argentinian Primera Division", "2025-11-10", "19:00:00", "Deportivo Riestra", "CA Independiente", 1.09, 1000, 1.09, 1000, 1.09, 1000, 0, 0, 0, 0, 1.25, 1.01, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  ,@("Brazilian Serie B", "2025-11-10", "19:00:00", "Botafogo SP", "Amazonas FC", 1.82, 2.48, 3.6, 4.9, 3.1, 4.1, 0, 0, 0, 0, 1.25, 1.01, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  ,@("Brazilian Serie B", "2025-11-10", "19:00:00", "Chapecoense", "America MG", 2.04, 2.22, 3.95, 4.7, 3.25, 3.7, 0, 0, 0, 0, 1.7, 2.2, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  ,@("Paraguayan Primera Division", "2025-11-10", "19:30:00", "Club Atletico Tembetary", "Nacional (Par)", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1.24, 1.01, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value = $row[$j]
  }
}

$ws.Range("B2:B12").Style = "Normal"

$ws.Range("A1").Select() | Out-Null